# Generate Report for Handback
# Refresh the handback-status report: the 5afef0f6... file has been
# re-handed-off/handed-back, so its recorded timestamps move forward on
# the Overview sheet and on each per-language detail sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview: "Latest HO Xliff Generate Date" for 5afef0f6-2833-4e05-be18-cb778151c15b.md
$wsOverview.Range("G2").Value = "2016-08-15 10:46:58"

# zh-cn detail sheet: Correspond Handoff / Handback datetimes for the same file
$wsZhCn.Range("H2").Value = "2016-08-15 10:46:54"
$wsZhCn.Range("K2").Value = "2016-08-15 10:47:14"

# de-de detail sheet: Correspond Handoff / Handback datetimes for the same file
$wsDeDe.Range("H2").Value = "2016-08-15 10:46:58"
$wsDeDe.Range("K2").Value = "2016-08-15 10:47:20"
